$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values, forcing text so Excel does not
# reinterpret numeric-looking strings (e.g. "1.002") as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.427.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5077'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08223'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.106'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.280'
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.497'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.819.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001144'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06634'
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.095'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.447.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.267'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.025.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '155.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.402'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.90'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.106'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1094'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.782'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.651'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07041'
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02336'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.210'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.827'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6269'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.406'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.739'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5891'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06886'
$ws.Range("D51").Style = "Normal"

# Update remaining columns (Coin name, Link, Volume(1h)).
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B28").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E7").Value = '  -5.07%  '
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("E9").Value = '  +7.85%  '
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("E17").Value = '  +6.46%  '
$ws.Range("E18").Value = '  +3.15%  '
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("E28").Value = '  -1.29%  '
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  -6.52%  '
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("E51").Value = '  -0.08%  '
